$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 20,12
$arr[0,0] = 0.7519746824633965
$arr[0,1] = 0.9987653331488298
$arr[0,2] = 0.8493545730501202
$arr[0,3] = 0.8476067057428658
$arr[0,4] = 0.7356418701799486
$arr[0,5] = 0.9157278457125603
$arr[0,6] = 0.9116010388328724
$arr[0,7] = 0.8227567499526767
$arr[0,8] = 0.9116560933846574
$arr[0,9] = 0.8874281520046929
$arr[0,10] = 0.9031170874463322
$arr[0,11] = 0.8970926901064841
$arr[1,0] = 0.7251677180060386
$arr[1,1] = 0.9987812031962048
$arr[1,2] = 0.8506635905784974
$arr[1,3] = 0.8566337270199571
$arr[1,4] = 0.7862917737789203
$arr[1,5] = 0.9187588512911203
$arr[1,6] = 0.9130835344815977
$arr[1,7] = 0.8375648955377621
$arr[1,8] = 0.9077065739138729
$arr[1,9] = 0.8862806288215519
$arr[1,10] = 0.9037489449928443
$arr[1,11] = 0.8920802300244219
$arr[2,0] = 0.7496042424335313
$arr[2,1] = 0.9988051795267712
$arr[2,2] = 0.8555739516448069
$arr[2,3] = 0.8430070090246918
$arr[2,4] = 0.8202594794344473
$arr[2,5] = 0.9220655049459398
$arr[2,6] = 0.9125323223578847
$arr[2,7] = 0.8152589912088948
$arr[2,8] = 0.9063210194461954
$arr[2,9] = 0.8912535818235149
$arr[2,10] = 0.9210166599390849
$arr[2,11] = 0.8980106669290066
$arr[3,0] = 0.7439582653666696
$arr[3,1] = 0.9988640927961631
$arr[3,2] = 0.850669296430419
$arr[3,3] = 0.8609344900385199
$arr[3,4] = 0.8083410989717225
$arr[3,5] = 0.9217723340090866
$arr[3,6] = 0.9141010399187907
$arr[3,7] = 0.8454143369292725
$arr[3,8] = 0.8998464675998574
$arr[3,9] = 0.8935239502718801
$arr[3,10] = 0.8858138233459322
$arr[3,11] = 0.8939257658909966
$arr[4,0] = 0.7384442074886308
$arr[4,1] = 0.9988178527300706
$arr[4,2] = 0.8537666230485985
$arr[4,3] = 0.8347851676086313
$arr[4,4] = 0.8189042416452442
$arr[4,5] = 0.9155719371836899
$arr[4,6] = 0.9073789885235811
$arr[4,7] = 0.8159416975510011
$arr[4,8] = 0.9052605244765146
$arr[4,9] = 0.8900454214706346
$arr[4,10] = 0.9385171645077245
$arr[4,11] = 0.8904311296992014
$arr[5,0] = 0.7404963445794899
$arr[5,1] = 0.9988906380192902
$arr[5,2] = 0.8438802836188796
$arr[5,3] = 0.8628756179580732
$arr[5,4] = 0.7774783097686375
$arr[5,5] = 0.9184964576863355
$arr[5,6] = 0.9105854880486006
$arr[5,7] = 0.8332394872351689
$arr[5,8] = 0.9088514539566727
$arr[5,9] = 0.8913445375781233
$arr[5,10] = 0.9035482642838795
$arr[5,11] = 0.9033866126144198
$arr[6,0] = 0.7432513208340358
$arr[6,1] = 0.9988342365559577
$arr[6,2] = 0.8506597866772162
$arr[6,3] = 0.8644393906100654
$arr[6,4] = 0.8164058483290489
$arr[6,5] = 0.9198239267023234
$arr[6,6] = 0.9080127304372907
$arr[6,7] = 0.8261284833157332
$arr[6,8] = 0.9104607313223978
$arr[6,9] = 0.8904972036393584
$arr[6,10] = 0.9182163039888445
$arr[6,11] = 0.8938583709329059
$arr[7,0] = 0.7463604535830565
$arr[7,1] = 0.9988123724259411
$arr[7,2] = 0.8422636255743891
$arr[7,3] = 0.8353791527520146
$arr[7,4] = 0.8075297236503856
$arr[7,5] = 0.9192050102800783
$arr[7,6] = 0.912920212370868
$arr[7,7] = 0.842354890070346
$arr[7,8] = 0.905559132141436
$arr[7,9] = 0.8882389978791092
$arr[7,10] = 0.9206267659902388
$arr[7,11] = 0.8948906955990132
$arr[8,0] = 0.7323625139308853
$arr[8,1] = 0.9988180239895746
$arr[8,2] = 0.8512603275919783
$arr[8,3] = 0.8530861216341812
$arr[8,4] = 0.7939644922879178
$arr[8,5] = 0.9157626667816885
$arr[8,6] = 0.9113738647267243
$arr[8,7] = 0.8258128337261736
$arr[8,8] = 0.9052828781382678
$arr[8,9] = 0.8932926829268293
$arr[8,10] = 0.8845397875307328
$arr[8,11] = 0.8947147021302556
$arr[9,0] = 0.7380922405242611
$arr[9,1] = 0.9987696146364309
$arr[9,2] = 0.8546301086394206
$arr[9,3] = 0.8621791412605373
$arr[9,4] = 0.8026582583547557
$arr[9,5] = 0.9101998324994249
$arr[9,6] = 0.9070766688718048
$arr[9,7] = 0.8432192030256239
$arr[9,8] = 0.9077464379439997
$arr[9,9] = 0.8864893334987928
$arr[9,10] = 0.9219546989101317
$arr[9,11] = 0.9020617107842792
$arr[10,0] = 0.7481227516262735
$arr[10,1] = 0.9989516064027306
$arr[10,2] = 0.8485319793980706
$arr[10,3] = 0.850761426837724
$arr[10,4] = 0.8056973007712084
$arr[10,5] = 0.9113583631671268
$arr[10,6] = 0.9061953375881061
$arr[10,7] = 0.8455032184727557
$arr[10,8] = 0.9027256192336868
$arr[10,9] = 0.8927800873175246
$arr[10,10] = 0.9214845326776999
$arr[10,11] = 0.8896444293116325
$arr[11,0] = 0.7441230275070667
$arr[11,1] = 0.9987949610430297
$arr[11,2] = 0.8506288799793067
$arr[11,3] = 0.8478870046144492
$arr[11,4] = 0.8022670308483291
$arr[11,5] = 0.916050895373246
$arr[11,6] = 0.9070480006289638
$arr[11,7] = 0.8445807721836323
$arr[11,8] = 0.9049382591862373
$arr[11,9] = 0.8891545486338305
$arr[11,10] = 0.8835547319364426
$arr[11,11] = 0.9052541875905121
$arr[12,0] = 0.7246145879632773
$arr[12,1] = 0.9988189944600975
$arr[12,2] = 0.8385386552448191
$arr[12,3] = 0.825984093621376
$arr[12,4] = 0.8076012210796916
$arr[12,5] = 0.9178748454393836
$arr[12,6] = 0.9135656822021033
$arr[12,7] = 0.8244637560283306
$arr[12,8] = 0.9072028714023644
$arr[12,9] = 0.8907896369666749
$arr[12,10] = 0.9219294704781476
$arr[12,11] = 0.9072393044840324
$arr[13,0] = 0.7461130366811978
$arr[13,1] = 0.9987930201019839
$arr[13,2] = 0.8444908097745047
$arr[13,3] = 0.8346003722306898
$arr[13,4] = 0.7874212724935732
$arr[13,5] = 0.9155899093484012
$arr[13,6] = 0.9015780130214635
$arr[13,7] = 0.8348239809130489
$arr[13,8] = 0.9084278520664471
$arr[13,9] = 0.8921370512849438
$arr[13,10] = 0.8847462019742394
$arr[13,11] = 0.8903813021472482
$arr[14,0] = 0.7497438249777215
$arr[14,1] = 0.9988485652677962
$arr[14,2] = 0.8507139922704725
$arr[14,3] = 0.8097593700185806
$arr[14,4] = 0.7779884318766068
$arr[14,5] = 0.9161133486456176
$arr[14,6] = 0.9124406708542572
$arr[14,7] = 0.8469613562266567
$arr[14,8] = 0.9078773931457711
$arr[14,9] = 0.8892624263892963
$arr[14,10] = 0.8847381747458809
$arr[14,11] = 0.8923232990438859
$arr[15,0] = 0.7298949135026133
$arr[15,1] = 0.9988384038705561
$arr[15,2] = 0.8517610160981102
$arr[15,3] = 0.8437562844016394
$arr[15,4] = 0.836577763496144
$arr[15,5] = 0.9181770024585921
$arr[15,6] = 0.9064616047526735
$arr[15,7] = 0.8313417461716076
$arr[15,8] = 0.8944717159117835
$arr[15,9] = 0.8935555027526455
$arr[15,10] = 0.9207070382738249
$arr[15,11] = 0.896176310320883
$arr[16,0] = 0.7333770874432638
$arr[16,1] = 0.9988459392887341
$arr[16,2] = 0.8488515071056876
$arr[16,3] = 0.8392489853724364
$arr[16,4] = 0.8449638496143959
$arr[16,5] = 0.9173354558459857
$arr[16,6] = 0.9125197257057275
$arr[16,7] = 0.8482609485246143
$arr[16,8] = 0.9170168740341356
$arr[16,9] = 0.890846396177884
$arr[16,10] = 0.9223847290007707
$arr[16,11] = 0.8905320624326449
$arr[17,0] = 0.736002607730355
$arr[17,1] = 0.9988644353151713
$arr[17,2] = 0.8610762948479961
$arr[17,3] = 0.8418081684215757
$arr[17,4] = 0.8083603791773779
$arr[17,5] = 0.9195711931360708
$arr[17,6] = 0.9096876508069025
$arr[17,7] = 0.8337521175427184
$arr[17,8] = 0.9043786725203455
$arr[17,9] = 0.8946062532433836
$arr[17,10] = 0.9192747513852703
$arr[17,11] = 0.8976337662155145
$arr[18,0] = 0.73438618713038
$arr[18,1] = 0.9988002700876552
$arr[18,2] = 0.843594991022793
$arr[18,3] = 0.8410472462771108
$arr[18,4] = 0.8072734575835476
$arr[18,5] = 0.9155566608436854
$arr[18,6] = 0.9167005111200205
$arr[18,7] = 0.8506338455152296
$arr[18,8] = 0.9046396515213155
$arr[18,9] = 0.8856541622481442
$arr[18,10] = 0.937540136141793
$arr[18,11] = 0.8910373649146319
$arr[19,0] = 0.7543547016874709
$arr[19,1] = 0.9988611242980929
$arr[19,2] = 0.8424409824716228
$arr[19,3] = 0.8304494642875168
$arr[19,4] = 0.8242753856041132
$arr[19,5] = 0.9190989745082816
$arr[19,6] = 0.9051202784815734
$arr[19,7] = 0.8293950001489366
$arr[19,8] = 0.9095406918607337
$arr[19,9] = 0.8857846026714199
$arr[19,10] = 0.9213813254559466
$arr[19,11] = 0.8934619863689705

$range = $ws.Range("B2:M21")
$range.Value = $arr
